# Updated cryptos list on Fri Mar  1 04:56:49 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) values for
# each coin row on the sheet, and fixes two rows where the ranking of a
# pair of coins swapped places (Toncoin / InjectiveProtocol at rows 34-35,
# and TheGraph / ARBITRUM at rows 42-43) by rewriting the Coin name, Link,
# Price and Volume(1h) cells for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "61.477.36"
$ws.Cells.Item(2, 5).Value = "  -1.72%  "
$ws.Cells.Item(3, 4).Value = "3.378.67"
$ws.Cells.Item(3, 5).Value = "  -2.01%  "
$ws.Cells.Item(4, 5).Value = "  -0.03%  "
$ws.Cells.Item(5, 4).Value = "407.25"
$ws.Cells.Item(5, 5).Value = "  -1.94%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "136.00"
$ws.Cells.Item(6, 5).Value = "  +10.82%  "
$ws.Cells.Item(7, 4).Value = "0.595"
$ws.Cells.Item(7, 5).Value = "  +0.46%  "
$ws.Cells.Item(8, 5).Value = "  -0.01%  "
$ws.Cells.Item(9, 4).Value = "0.675"
$ws.Cells.Item(9, 5).Value = "  +3.30%  "
$ws.Cells.Item(10, 4).Value = "0.121"
$ws.Cells.Item(10, 5).Value = "  -5.74%  "
$ws.Cells.Item(11, 4).Value = "43.14"
$ws.Cells.Item(11, 5).Value = "  +4.79%  "
$ws.Cells.Item(12, 5).Value = "  -1.15%  "
$ws.Cells.Item(13, 4).Value = "3.901.06"
$ws.Cells.Item(13, 5).Value = "  -2.17%  "
$ws.Cells.Item(14, 4).Value = "8.41"
$ws.Cells.Item(14, 5).Value = "  -0.91%  "
$ws.Cells.Item(15, 4).Value = "19.74"
$ws.Cells.Item(15, 5).Value = "  +0.56%  "
$ws.Cells.Item(16, 4).Value = "3.368.89"
$ws.Cells.Item(16, 5).Value = "  -1.83%  "
$ws.Cells.Item(17, 4).Value = "61.391.52"
$ws.Cells.Item(17, 5).Value = "  -1.55%  "
$ws.Cells.Item(18, 5).Value = "  -0.53%  "
$ws.Cells.Item(19, 4).Value = "11.05"
$ws.Cells.Item(20, 4).Value = "0.0000128"
$ws.Cells.Item(20, 5).Value = "  -4.60%  "
$ws.Cells.Item(21, 4).Value = "3.22"
$ws.Cells.Item(21, 5).Value = "  -2.71%  "
$ws.Cells.Item(22, 4).Value = "83.53"
$ws.Cells.Item(22, 5).Value = "  +1.83%  "
$ws.Cells.Item(23, 4).Value = "314.92"
$ws.Cells.Item(23, 5).Value = "  -0.71%  "
$ws.Cells.Item(24, 4).Value = "12.89"
$ws.Cells.Item(24, 5).Value = "  -0.38%  "
$ws.Cells.Item(25, 5).Value = "  -0.52%  "
$ws.Cells.Item(26, 4).Value = "4.79"
$ws.Cells.Item(26, 5).Value = "  +11.43%  "
$ws.Cells.Item(27, 4).Value = "8.34"
$ws.Cells.Item(27, 5).Value = "  +7.54%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "29.50"
$ws.Cells.Item(28, 5).Value = "  -5.06%  "
$ws.Cells.Item(29, 4).Value = "7.77"
$ws.Cells.Item(29, 5).Value = "  -1.29%  "
$ws.Cells.Item(30, 5).Value = "  +1.62%  "
$ws.Cells.Item(31, 4).Value = "0.173"
$ws.Cells.Item(31, 5).Value = "  -1.18%  "
$ws.Cells.Item(32, 4).Value = "11.37"
$ws.Cells.Item(32, 5).Value = "  -0.18%  "
$ws.Cells.Item(33, 5).Value = "  -0.07%  "
$ws.Cells.Item(34, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(34, 4).Value = "41.13"
$ws.Cells.Item(34, 5).Value = "  -1.87%  "
$ws.Cells.Item(35, 2).Value = "Toncoin"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(35, 4).Value = "2.49"
$ws.Cells.Item(35, 5).Value = "  -2.62%  "
$ws.Cells.Item(36, 5).Value = "  +0.14%  "
$ws.Cells.Item(37, 4).Value = "52.07"
$ws.Cells.Item(37, 5).Value = "  -1.22%  "
$ws.Cells.Item(38, 4).Value = "0.997"
$ws.Cells.Item(38, 5).Value = "  -0.05%  "
$ws.Cells.Item(39, 4).Value = "3.43"
$ws.Cells.Item(39, 5).Value = "  -1.98%  "
$ws.Cells.Item(40, 4).Value = "2.95"
$ws.Cells.Item(40, 5).Value = "  -3.26%  "
$ws.Cells.Item(41, 4).Value = "138.24"
$ws.Cells.Item(41, 5).Value = "  +2.44%  "
$ws.Cells.Item(42, 2).Value = "ARBITRUM"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(42, 4).Value = "1.98"
$ws.Cells.Item(42, 5).Value = "  -0.46%  "
$ws.Cells.Item(43, 2).Value = "TheGraph"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(43, 4).Value = "0.298"
$ws.Cells.Item(43, 5).Value = "  +6.05%  "
$ws.Cells.Item(44, 5).Value = "  -0.50%  "
$ws.Cells.Item(45, 4).Value = "4.04"
$ws.Cells.Item(45, 5).Value = "  +4.41%  "
$ws.Cells.Item(46, 4).Value = "16.73"
$ws.Cells.Item(46, 5).Value = "  -1.93%  "
$ws.Cells.Item(47, 5).Value = "  +1.20%  "
$ws.Cells.Item(48, 4).Value = "21.37"
$ws.Cells.Item(48, 5).Value = "  -2.34%  "
$ws.Cells.Item(49, 4).Value = "2.131.64"
$ws.Cells.Item(49, 5).Value = "  -3.22%  "
$ws.Cells.Item(51, 5).Value = "  +1.04%  "
